# Fruta / hortaliza, semanal
# Insert a new weekly record as row 5, shifting existing rows 5-17 down to 6-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes old rows 5..17 down to 6..18,
# carrying their formatting/styles with them automatically).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44804
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100112017
$ws.Range("G5").Value = "Ramas de apio"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 5500
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5750
$ws.Range("N5").Value = "$/atado 7 kilos"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 5750
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
